$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data: date 2023-09-10 (serial 45179) with 8 hours,
# matching the formatting of the row above it
$ws.Range("A24").Copy()
$ws.Range("A25").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A25").Value2 = 45179
$ws.Range("B25").Value2 = 8

# Move the active selection to B26, matching the author's next input cell
$ws.Range("B26").Select()
